# Excel COM-interop script: daily attendance rollover.
# 1. Archive the current "Sheet1" attendance snapshot into a new sheet
#    named after today's date, refreshing each person's check-in time.
# 2. Reset "Sheet1" with the next day's roll-call, storing Time as a
#    real Excel date/time serial instead of free-text.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Archive current Sheet1 (style untouched by Copy) into a dated sheet ---
$ws1.Copy($null, $ws1)
$archive = $wb.Worksheets.Item(2)
$archive.Name = "2024-08-21"

# Re-stamp the archive header format onto its own style slot so the header
# restyle below (step 2) only touches the live "Sheet1", not this archive.
$archive.Range("A1:B1").NumberFormat = "General"

$archive.Range("A1").Select()

# Refresh the archived check-in times to reflect when this rollover ran.
$archive.Range("B2").Value = "2024-08-21 00:08:42"
$archive.Range("B3").Value = "2024-08-21 00:08:42"
$archive.Range("B4").Value = "2024-08-21 00:08:42"
$archive.Range("B5").Value = "2024-08-21 00:08:49"
$archive.Range("B6").Value = "2024-08-21 00:08:49"
$archive.Range("B7").Value = "2024-08-21 00:08:51"

# --- 2. Reset Sheet1 with the new day's attendance, Time as a date serial ---
$ws1.Range("B2").Value = 45524.98214120371
$ws1.Range("B3").Value = 45525.98214120371
$ws1.Range("B4").Value = 45526.98214120371
$ws1.Range("B5").Value = 45527.98214120371
$ws1.Range("B6").Value = 45528.98214120371
$ws1.Range("B7").Value = 45529.98214120371
$ws1.Range("B2:B7").NumberFormat = "m/d/yyyy h:mm"

# Refresh the header look on the reset sheet: plain Calibri 11 bold with an
# explicit (auto-color) thin box border.
$header = $ws1.Range("A1:B1")
$header.Font.Name = "Calibri"
$header.Font.Bold = $true
$header.Font.Size = 11
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$ws1.Columns.Item(1).ColumnWidth = 15
$ws1.Columns.Item(2).ColumnWidth = 24.26953125

$ws1.Range("D10").Select()
$ws1.Activate()
